$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: agent_minimax vs agent_random -> agent_random vs agent_random ---
$ws.Range("A2").Value = "agent_random"
$ws.Range("D2").Value = "agent_random"
$ws.Range("E2").Value = "'0.029"

$ws.Range("G2").Value = "2 1 2 0 0 0
1 2 0 0 0 0
2 1 1 0 0 0
2 1 1 0 0 0
2 1 2 0 0 0
1 1 0 0 0 0
1 2 2 0 0 0"

# --- Row 3: winner column + match time + final board ---
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = "'0.0426"

$ws.Range("G3").Value = "1 2 2 0 0 0
0 0 0 0 0 0
2 1 1 1 1 0
1 2 1 2 1 2
2 1 0 0 0 0
2 2 0 0 0 0
1 1 2 1 2 0"

# Re-fit the rows so the embedded newlines in G2/G3 don't leave a stray
# custom row height behind (keeps row sizing identical to the source file).
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()

# --- Selection / active view state, matches the recorded sheetView change ---
$ws.Range("H9").Select()
